$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62 - new/updated record (week shift: newest entry)
$ws.Range("D62").Value = 44588
$ws.Range("K62").Value = "Lapins"
$ws.Range("Q62").Value = "$/bandeja 10 kilos"
$ws.Range("S62").Value = 1050
$ws.Range("T62").Value = 10

# Row 63 - shifted from old row 62
$ws.Range("D63").Value = 44536
$ws.Range("K63").Value = "Santina"
$ws.Range("M63").Value = 400
$ws.Range("N63").Value = 10000
$ws.Range("O63").Value = 11000
$ws.Range("P63").Value = 10500
$ws.Range("Q63").Value = "$/bandeja 8 kilos"
$ws.Range("S63").Value = 1312

# Row 64 - shifted from old row 63
$ws.Range("K64").Value = "Lapins"
$ws.Range("N64").Value = 9000
$ws.Range("O64").Value = 10000
$ws.Range("P64").Value = 9500
$ws.Range("Q64").Value = "$/caja 8 kilos"
$ws.Range("S64").Value = 1188
$ws.Range("T64").Value = 8

# Row 65 - shifted from old row 64
$ws.Range("D65").Value = 44187
$ws.Range("K65").Value = "Rainier"
$ws.Range("N65").Value = 16000
$ws.Range("O65").Value = 17000
$ws.Range("P65").Value = 16500
$ws.Range("Q65").Value = "$/bandeja 10 kilos"
$ws.Range("S65").Value = 1650
$ws.Range("T65").Value = 10

# Row 66 - shifted from old row 65 (old row 66 data is dropped)
$ws.Range("D66").Value = 44533
$ws.Range("M66").Value = 600
$ws.Range("N66").Value = 10000
$ws.Range("O66").Value = 11000
$ws.Range("P66").Value = 10500
$ws.Range("Q66").Value = "$/caja 8 kilos"
$ws.Range("S66").Value = 1312
$ws.Range("T66").Value = 8
